# Add a new "UK" Test Data sheet, cloned from the existing "Poland" sheet,
# matching the layout/styles used by every other country sheet in this
# workbook, then fill in the UK-specific values.

$wb = $excel.ActiveWorkbook

# Poland is the right-most existing country sheet - use it as the template
# for the new UK sheet (same columns/styles/merged cells/page setup).
$template = $wb.Worksheets.Item("Poland")

# Copy it to the end of the workbook (after itself); the copy becomes the
# active sheet/tab automatically, just like it does in real Excel.
$template.Copy([System.Reflection.Missing]::Value, $template)

$ukSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ukSheet.Name = "UK"

# Fill in the UK-specific values (mirrors how every other country sheet
# stores its Jira reference in B4 and its "<Country> Market" label in B2).
$ukSheet.Range("B4").Value = "NGC-2741/T3340"
$ukSheet.Range("B2").Value = "UK Market"

# Leave the selection on B4, same as the authored sheet.
[void]$ukSheet.Range("B4").Select()
